# add spu - BNIA nb
# Append a new data row (row 16) documenting the BNIA-JFI neighborhood
# shapefile source, and move the active selection down to B17 to match
# where Excel would leave the cursor after the new last row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the url column (C) first, then filename (A) and remarks (B), so the
# new shared strings land in the same append order as the source edit:
#   .. , <geojson url>, <Vital Signs 15 Census Demographics>, <remarks>
$ws.Range("C16").Value = "https://opendata.arcgis.com/datasets/794586676bcc4f5fb629c08c51474cf6_0.geojson"
$ws.Range("A16").Value = "Vital Signs 15 Census Demographics"
$ws.Range("B16").Value = "extracted the shape of neighborhood defined by BNIA-JFI"

# Move the selection to B17 (one row below the new last data row).
$ws.Range("B17").Select() | Out-Null
